$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.963.10"
$ws.Range("E2").Value = "  -3.72%  "

$ws.Range("D3").Value = "2.518.58"
$ws.Range("E3").Value = "  -4.75%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'577.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "

$ws.Range("D6").Value = "'168.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.57%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").Value = "2.517.91"
$ws.Range("E9").Value = "  -4.71%  "

$ws.Range("E10").Value = "  -6.32%  "

$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("E12").Value = "  -3.71%  "

$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("D14").Value = "2.979.17"
$ws.Range("E14").Value = "  -4.69%  "

$ws.Range("D15").Value = "69.858.47"
$ws.Range("E15").Value = "  -3.60%  "

$ws.Range("E16").Value = "  -6.14%  "

$ws.Range("D18").Value = "2.522.89"
$ws.Range("E18").Value = "  -4.33%  "

$ws.Range("E19").Value = "  -1.11%  "

$ws.Range("D20").Value = "'11.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.34%  "

$ws.Range("D21").Value = "'350.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.00%  "

$ws.Range("E22").Value = "  -4.11%  "

$ws.Range("E23").Value = "  -4.12%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "'69.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.53%  "

$ws.Range("E26").Value = "  -5.84%  "

$ws.Range("E27").Value = "  -5.47%  "

$ws.Range("D28").Value = "2.648.46"
$ws.Range("E28").Value = "  -4.67%  "

$ws.Range("D29").Value = "'1.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.97%  "

$ws.Range("E30").Value = "  -4.63%  "

$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D32").Value = "'1.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "

$ws.Range("D33").Value = "'467.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.16%  "

$ws.Range("E34").Value = "  -2.02%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("E36").Value = "  +4.78%  "

$ws.Range("D37").Value = "'152.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.53%  "

$ws.Range("D38").Value = "'19.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Value = "'18.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.81%  "

$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("E41").Value = "  -2.51%  "

$ws.Range("D42").Value = "'0.321"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("E43").Value = "  -7.06%  "

$ws.Range("E44").Value = "  -13.92%  "

$ws.Range("E45").Value = "  -10.01%  "

$ws.Range("D46").Value = "'38.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.39%  "

$ws.Range("D47").Value = "'143.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.76%  "

$ws.Range("D48").Value = "'0.533"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("D49").Value = "'3.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.54%  "

$ws.Range("E50").Value = "  -4.45%  "

$ws.Range("E51").Value = "  -1.23%  "
